# Apply updated dSF (column F) values for a set of rows.
# These correspond to a "repull" of source data where the dSF value
# (distance-to-set final) was recalculated for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -7
    5  = -4
    8  = -5
    9  = 2
    17 = -3
    20 = 1
    21 = 2
    26 = 1
    31 = -2
    44 = -7
    46 = -3
    54 = -6
    56 = -1
    58 = -4
    63 = 3
    68 = 3
    69 = 3
    72 = -1
    79 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
